$d = $word.ActiveDocument
$q1 = [char]0x2018
$q2 = [char]0x2019

# ---------------------------------------------------------------------
# 1) Remove the "Inspect code for aliasing violations." list item
#    entirely (whole paragraph, including its paragraph mark).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Inspect code for aliasing violations.*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) The next item - "Use 'auto' by default where appropriate." - used
#    to be split across two runs (with a _GoBack bookmark sitting in
#    between them). Merge the two runs into a single run while leaving
#    the bookmark in place immediately before the (now single) run.
# ---------------------------------------------------------------------
$rngFirst = $d.Content
$rngFirst.Find.Execute("Use " + $q1 + "auto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$firstStart = $rngFirst.Start
$firstEnd = $rngFirst.End

$rngSecond = $d.Content
$rngSecond.Find.Execute($q2 + " by default where appropriate.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Rewrite the second run's text so it carries the full merged sentence.
$rngSecond.Text = "Use " + $q1 + "auto" + $q2 + " by default where appropriate."

# Delete the now-redundant first run's original text (leaves the
# bookmark, which sits between the two original runs, untouched).
$rngDup = $d.Range($firstStart, $firstEnd)
$rngDup.Delete()

# ---------------------------------------------------------------------
# 3) "Assertions where appropriate (for 'impossible' or rare cases)."
#    loses its cached <w:lastRenderedPageBreak/> marker (it shifts to
#    the following item - see step 4).
# ---------------------------------------------------------------------
$rngAssert = $d.Content
$rngAssert.Find.Execute("Assertions where appropriate (for " + $q1 + "impossible" + $q2 + " or rare cases).", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngAssert.Delete()

$assertXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009F2B2C" w:rsidRDefault="009F2B2C" w:rsidP="00E124A1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:t>Assertions where appropriate</w:t></w:r><w:r w:rsidR="0009236D" w:rsidRPr="008B06FC"><w:t xml:space="preserve"> (for ' + $q1 + 'impossible' + $q2 + ' </w:t></w:r><w:r w:rsidR="00BD6B7B"><w:t xml:space="preserve">or rare </w:t></w:r><w:r w:rsidR="0009236D" w:rsidRPr="008B06FC"><w:t>cases).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngAssert.InsertXML($assertXml)

# ---------------------------------------------------------------------
# 4) "Especially check for things like integer overflows ..." gains the
#    <w:lastRenderedPageBreak/> marker that used to sit on the previous
#    item.
# ---------------------------------------------------------------------
$rngEsp = $d.Content
$rngEsp.Find.Execute("Especially check for things like integer overflows that would otherwise be hard to debug.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngEsp.Delete()

$espXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BD6B7B" w:rsidRPr="00277316" w:rsidRDefault="00BD6B7B" w:rsidP="00BD6B7B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="00277316"><w:lastRenderedPageBreak/><w:t>Especially check fo</w:t></w:r><w:r w:rsidR="009C66EE" w:rsidRPr="00277316"><w:t>r things like integer overflows that would otherwise be hard to debug.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngEsp.InsertXML($espXml)

Write-Host "Done."
